# AFDP-4158 Document Level - ACL Implementation
# File/Folder Assignment Rules
#
# Inserts a new "File - default access" rule row above the existing
# "Folder - default access" row on Sheet1 (new row 25), pushing the
# remaining rule rows (old 25-35) down to rows 26-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25; everything below (old rows 25-35,
# "Folder - default access" ... "Person - Default group") shifts down
# to rows 26-36.
$ws.Rows.Item(25).Insert()

# The freshly inserted row doesn't carry the table's formatting, so copy
# it from the row directly below (the old row 25, now at row 26, which
# has the same cell layout/style pattern we want for the new rule row).
$ws.Range("B26:H26").Copy()
$ws.Range("B25:H25").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Populate the new "File - default access" rule row.
$ws.Range("B25").Value = "File - default access"
$ws.Range("C25").Value = "FILE"
$ws.Range("D25").Value = "participants.?[participantType == '*'].isEmpty()"
$ws.Range("G25").Value = "*, *"

# Match the author's recorded selection for the sheet.
$ws.Range("D25").Select() | Out-Null

Write-Host "Edit applied successfully"
